$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (rows 3-5, columns E/F) ---
$wsSchedule.Range("E3").Value = -205.23275175
$wsSchedule.Range("F3").Value = -6.032708752204585
$wsSchedule.Range("E4").Value = 630.9916275
$wsSchedule.Range("F4").Value = 30.35072763347764
$wsSchedule.Range("E5").Value = -13.48967100000001
$wsSchedule.Range("F5").Value = -0.3756522138680036

# --- Detailed sheet updates (rows 32-97, column B, and C34) ---
$wsDetailed.Range("B32").Value = -24.08764
$wsDetailed.Range("B33").Value = -5.83482
$wsDetailed.Range("B34").Value = -5.21462
$wsDetailed.Range("C34").Value = "historical"
$wsDetailed.Range("B35").Value = -4.38094
$wsDetailed.Range("B37").Value = 47.82322
$wsDetailed.Range("B38").Value = 57.09
$wsDetailed.Range("B39").Value = 65
$wsDetailed.Range("B40").Value = 77.94
$wsDetailed.Range("B41").Value = 77.94
$wsDetailed.Range("B42").Value = 74.18452000000001
$wsDetailed.Range("B45").Value = 58.64889
$wsDetailed.Range("B49").Value = 60.9348
$wsDetailed.Range("B51").Value = 62.02711
$wsDetailed.Range("B56").Value = 63.71733
$wsDetailed.Range("B57").Value = 63.33647
$wsDetailed.Range("B58").Value = 63.32506
$wsDetailed.Range("B60").Value = 62.3992
$wsDetailed.Range("B62").Value = 63.14049
$wsDetailed.Range("B63").Value = 56.98
$wsDetailed.Range("B66").Value = -5.50985
$wsDetailed.Range("B67").Value = -7.61395
$wsDetailed.Range("B68").Value = -12.01
$wsDetailed.Range("B69").Value = -10
$wsDetailed.Range("B71").Value = -7.85417
$wsDetailed.Range("B72").Value = -7.70805
$wsDetailed.Range("B74").Value = -7.72886
$wsDetailed.Range("B75").Value = -7.9087
$wsDetailed.Range("B76").Value = -7.56872
$wsDetailed.Range("B77").Value = -5.87097
$wsDetailed.Range("B78").Value = -5.74313
$wsDetailed.Range("B79").Value = -5.50985
$wsDetailed.Range("B80").Value = -5.50985
$wsDetailed.Range("B84").Value = -11
$wsDetailed.Range("B85").Value = -8.256209999999999
$wsDetailed.Range("B86").Value = -5.97569
$wsDetailed.Range("B87").Value = -3.01569
$wsDetailed.Range("B88").Value = 10.11243
$wsDetailed.Range("B89").Value = 57.46182
$wsDetailed.Range("B90").Value = 36.61149
$wsDetailed.Range("B91").Value = 43.25197
$wsDetailed.Range("B92").Value = 38.26091
$wsDetailed.Range("B93").Value = 43.57572
$wsDetailed.Range("B94").Value = 30.67165
$wsDetailed.Range("B95").Value = 56.54508
$wsDetailed.Range("B96").Value = 56.54508
$wsDetailed.Range("B97").Value = 47.5613
